$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data rows after the TPM update: the "Sending cluster = ECs" rows have
# been removed, remaining rows (previously Sending=FAPs / Sending=MuSCs) keep
# their relative order and shift up, and several derived-specificity values
# were recomputed against the new TPM numbers.

$data = @(
    @("FAPs","Cxcl14","Cxcr4","ECs",3,1,574.5849303333333,1723.754791,0.9998767569401251,0.9998767569401251,3,1,13.71977066666667,41.159312,0.5515038136402627,0.5515038136402626,7883.173472695977,70948.5612542638,0.551435844622737,0.5514358446227369),
    @("FAPs","Cxcl14","Cxcr4","FAPs",3,1,574.5849303333333,1723.754791,0.9998767569401251,0.9998767569401251,3,1,4.280784333333333,12.842353,0.172077867958883,0.1720778679588829,2459.674167940358,22137.06751146322,0.172056660555899,0.1720566605558989),
    @("FAPs","Cxcl14","Cxcr4","MuSCs",3,1,574.5849303333333,1723.754791,0.9998767569401251,0.9998767569401251,3,1,6.876463666666666,20.629391,0.2764183184008545,0.2764183184008545,3951.112396851364,35560.01157166228,0.2763842517614893,0.2763842517614893),
    @("MuSCs","Cxcl14","Cxcr4","ECs",1,0.3333333333333333,0.07082233333333333,0.212467,0.0001232430598748645,0.0001232430598748645,3,1,13.71977066666667,41.159312,0.5515038136402627,0.5515038136402626,0.9716661714115556,8.744995542704,0.00006796901752568301,0.00006796901752568298),
    @("MuSCs","Cxcl14","Cxcr4","FAPs",1,0.3333333333333333,0.07082233333333333,0.212467,0.0001232430598748645,0.0001232430598748645,3,1,4.280784333333333,12.842353,0.172077867958883,0.1720778679588829,0.3031751349834444,2.728576214851,0.00002120740298399565,0.00002120740298399564),
    @("MuSCs","Cxcl14","Cxcr4","MuSCs",1,0.3333333333333333,0.07082233333333333,0.212467,0.0001232430598748645,0.0001232430598748645,3,1,6.876463666666666,20.629391,0.2764183184008545,0.2764183184008545,0.4870072019552222,4.383064817596999,0.00003406663936518587,0.00003406663936518586)
)

# Remove the old rows 8-10 (Sending=MuSCs) so the sheet only keeps 6 data rows.
$ws.Rows.Item(8).Resize(3).EntireRow.Delete() | Out-Null

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
